$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Remove the "NO_LABEL" labels that are no longer needed
$ws.Range("C3").ClearContents()
$ws.Range("C11").ClearContents()

# Drop the now-redundant standalone conditional-formatting blocks that only
# applied to C26 / G2:G9999 (their rules duplicate the patterns already
# covering the rest of the sheet).
$ws.Range("C26").FormatConditions.Delete()
$ws.Range("G2:G9999").FormatConditions.Delete()

# Update the active cell / selection on the frozen "survey" sheet
$ws.Range("A2").Select()
